# Update the "取得日時" (retrieved-at) timestamp column for all existing
# data rows on the "ランサーズ" sheet to reflect the latest scrape run.
#
# Commit message: Append: 2025-09-28 01:23 JST
# The rows scraped in this run happen to carry the same job listings as
# the rows already present, so every existing A2:A11 value is refreshed
# from "2025-09-27 18:28:10" to "2025-09-28 01:23:32".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-28 01:23:32"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 1
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}
